# Apply "Batterywise analysis" relabeling + value corrections to the
# "Analysis Results" sheet, and append the new "Time spent in 80-90 km/h" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap Starting/Ending SoC (%) values (rows 6 & 7) ---
$ws.Cells.Item(6, 2).Value = 96
$ws.Cells.Item(7, 2).Value = 9

# --- Relabel rows 8-30 (append units, no value change unless noted) ---
$ws.Cells.Item(8, 1).Value  = "Total distance covered (km)"
$ws.Cells.Item(9, 1).Value  = "Total energy consumption(WH/KM)"
$ws.Cells.Item(10, 1).Value = "Total SOC consumed(%)"
$ws.Cells.Item(12, 1).Value = "Peak Power(kW)"
$ws.Cells.Item(13, 1).Value = "Average Power(kW)"
$ws.Cells.Item(14, 1).Value = "Total Energy Regenerated(kWh)"

$ws.Cells.Item(15, 1).Value = "Regenerative Effectiveness(%)"
$ws.Cells.Item(15, 2).Value = 6.094035577472369

# Rows 16 & 17: label swap (Lowest<->Highest) plus corresponding value swap.
$ws.Cells.Item(16, 1).Value = "Highest Cell Voltage(V)"
$ws.Cells.Item(16, 2).Value = 3.337
$ws.Cells.Item(17, 1).Value = "Lowest Cell Voltage(V)"
$ws.Cells.Item(17, 2).Value = 2.921

$ws.Cells.Item(18, 1).Value = "Difference in Cell Voltage(V)"
$ws.Cells.Item(19, 1).Value = "Minimum Temperature(C)"
$ws.Cells.Item(20, 1).Value = "Maximum Temperature(C)"

$ws.Cells.Item(21, 1).Value = "Difference in Temperature(C)"
$ws.Cells.Item(21, 2).Value = 24

$ws.Cells.Item(22, 1).Value = "Maximum Fet Temperature-BMS(C)"
$ws.Cells.Item(23, 1).Value = "Maximum Afe Temperature-BMS(C)"
$ws.Cells.Item(24, 1).Value = "Maximum PCB Temperature-BMS(C)"
$ws.Cells.Item(25, 1).Value = "Maximum MCU Temperature(C)"
$ws.Cells.Item(26, 1).Value = "Maximum Motor Temperature(C)"
$ws.Cells.Item(27, 1).Value = "Abnormal Motor Temperature Detected(C)"

# Rows 28 & 29: label swap (lowest<->highest cell temp), values stay put.
$ws.Cells.Item(28, 1).Value = "highest cell temp(C)"
$ws.Cells.Item(29, 1).Value = "lowest cell temp(C)"

$ws.Cells.Item(30, 1).Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"

# --- Row 31: was "Maximum BMS Temperature in C" -> "Battery Voltage(V)" ---
$ws.Cells.Item(31, 1).Value = "Battery Voltage(V)"
$ws.Cells.Item(31, 2).Value = 55

# --- Row 32: was "Battery Voltage" -> "Total energy charged(kWh)" ---
$ws.Cells.Item(32, 1).Value = "Total energy charged(kWh)"
$ws.Cells.Item(32, 2).Value = 1.889129229166667

# --- Row 33: was "Total energy charged in kWh" -> "Electricity consumption units(kW)" ---
$ws.Cells.Item(33, 1).Value = "Electricity consumption units(kW)"
$ws.Cells.Item(33, 2).Value = [double]"1.177379670659553e-07"

# --- Row 34: was "Electricity consumption units in kW" -> "Idling time percentage" ---
$ws.Cells.Item(34, 1).Value = "Idling time percentage"
$ws.Cells.Item(34, 2).Value = 13.48716282092948

# --- Row 35: was "Idling time percentage" -> "Time spent in 0-10 km/h" ---
$ws.Cells.Item(35, 1).Value = "Time spent in 0-10 km/h"
$ws.Cells.Item(35, 2).Value = 13.26431124936162

# --- Row 36: was "Time spent in 0-10 km/h" -> "Time spent in 10-20 km/h" ---
$ws.Cells.Item(36, 1).Value = "Time spent in 10-20 km/h"
$ws.Cells.Item(36, 2).Value = 5.148799851432286

# --- Row 37: was "Time spent in 10-20 km/h" -> "Time spent in 20-30 km/h" ---
$ws.Cells.Item(37, 1).Value = "Time spent in 20-30 km/h"
$ws.Cells.Item(37, 2).Value = 8.284971447142393

# --- Row 38: was "Time spent in 20-30 km/h" -> "Time spent in 30-40 km/h" ---
$ws.Cells.Item(38, 1).Value = "Time spent in 30-40 km/h"
$ws.Cells.Item(38, 2).Value = 21.60731695993315

# --- Row 39: was "Time spent in 30-40 km/h" -> "Time spent in 40-50 km/h" ---
$ws.Cells.Item(39, 1).Value = "Time spent in 40-50 km/h"
$ws.Cells.Item(39, 2).Value = 14.78016621013046

# --- Row 40: was "Time spent in 40-50 km/h" -> "Time spent in 50-60 km/h" ---
$ws.Cells.Item(40, 1).Value = "Time spent in 50-60 km/h"
$ws.Cells.Item(40, 2).Value = 9.814754631134221

# --- Row 41: was "Time spent in 50-60 km/h" -> "Time spent in 60-70 km/h" ---
$ws.Cells.Item(41, 1).Value = "Time spent in 60-70 km/h"
$ws.Cells.Item(41, 2).Value = 10.51812990389526

# --- Row 42: was "Time spent in 60-70 km/h" -> "Time spent in 70-80 km/h" ---
$ws.Cells.Item(42, 1).Value = "Time spent in 70-80 km/h"
$ws.Cells.Item(42, 2).Value = 2.414225358651748

# --- New row 43: "Time spent in 80-90 km/h" ---
$ws.Cells.Item(43, 1).Value = "Time spent in 80-90 km/h"
$ws.Cells.Item(43, 2).Value = 0.03714192859464228
